$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. In-place numeric updates ---

# Estados Unidos (row 4): Casos totales, Nuevos casos, Recuperados updated
$ws.Cells.Item(4,2).Value = 215215
$ws.Cells.Item(4,3).Value = 212
$ws.Cells.Item(4,5).Value = 201227

# Turquia (row 13): Casos criticos updated
$ws.Cells.Item(13,6).Value = 979

# --- 2. Kazajistan: move from its old spot (row 79) up to right after
#        Bulgaria (new row 77), with refreshed case numbers ---
$ws.Rows(79).Delete()
$ws.Rows(77).Insert()
$ws.Cells.Item(77,1).Value = "Kazajistan"
$ws.Cells.Item(77,2).Value = 402
$ws.Cells.Item(77,3).Value = 22
$ws.Cells.Item(77,4).Value = 26
$ws.Cells.Item(77,5).Value = 373
$ws.Cells.Item(77,6).Value = 6
$ws.Cells.Item(77,7).Value = 0
$ws.Cells.Item(77,8).Value = 3

# --- 3. Honduras: move from its old spot (row 104) up to right after
#        Vietnam (new row 95), with refreshed case numbers ---
$ws.Rows(104).Delete()
$ws.Rows(95).Insert()
$ws.Cells.Item(95,1).Value = "Honduras"
$ws.Cells.Item(95,2).Value = 219
$ws.Cells.Item(95,3).Value = 47
$ws.Cells.Item(95,4).Value = 3
$ws.Cells.Item(95,5).Value = 202
$ws.Cells.Item(95,6).Value = 4
$ws.Cells.Item(95,7).Value = 4
$ws.Cells.Item(95,8).Value = 14

# --- 4. El Salvador: move from its old spot (row 143) up to right after
#        Jamaica (new row 137), with refreshed case numbers ---
$ws.Rows(143).Delete()
$ws.Rows(137).Insert()
$ws.Cells.Item(137,1).Value = "El Salvador"
$ws.Cells.Item(137,2).Value = 41
$ws.Cells.Item(137,3).Value = 8
$ws.Cells.Item(137,4).Value = 0
$ws.Cells.Item(137,5).Value = 39
$ws.Cells.Item(137,6).Value = 4
$ws.Cells.Item(137,7).Value = 0
$ws.Cells.Item(137,8).Value = 2

# --- 5. Update "last refreshed" timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 05:20"
